# "cropped high-level concept figure"
#
# The whole diagram on slide 2 ("High-level Representation") is shifted
# down (and the title box nudged slightly) to match a re-cropped figure.
# PowerPoint's Shape.Top / Shape.Left are expressed in points (a 32-bit
# float) while the OOXML stores EMU (1 pt = 12700 EMU). A naive
# "EMU / 12700.0" assignment loses sub-EMU precision when PowerPoint's
# single-precision float gets converted back to EMU on save (it truncates
# rather than rounds), so we hunt for the nearest representable float32
# that truncates back to the exact target EMU value.

function ConvertTo-PtForEmu($Emu) {
    $f = [float]($Emu / 12700.0)
    $best = $f
    $bestDiff = [Math]::Abs(([Math]::Floor([double]$f * 12700.0)) - $Emu)

    for ($i = 0; $i -lt 400; $i++) {
        $emuTest = [Math]::Floor([double]$f * 12700.0)
        $diff = [Math]::Abs($emuTest - $Emu)
        if ($diff -lt $bestDiff) {
            $bestDiff = $diff
            $best = $f
        }
        if ($emuTest -eq $Emu) {
            return $f
        }

        $stepFrac = 0.0000003 / ([Math]::Floor($i / 20) + 1)
        $step = [float]([Math]::Abs($f) * $stepFrac)
        if ($step -eq 0) {
            $step = [float]0.00000012
        }

        if ($emuTest -lt $Emu) {
            $nf = [float]($f + $step)
        } else {
            $nf = [float]($f - $step)
        }

        if ($nf -eq $f) {
            if ($emuTest -lt $Emu) {
                $nf = [float]($f * [float]1.0000001)
                if ($nf -eq $f) { $nf = [float]($f + [float]0.000001) }
            } else {
                $nf = [float]($f * [float]0.9999999)
                if ($nf -eq $f) { $nf = [float]($f - [float]0.000001) }
            }
        }
        $f = $nf
    }

    return $best
}

function Set-ShapeTopEmu($Shape, $Emu) {
    $Shape.Top = ConvertTo-PtForEmu $Emu
}

function Set-ShapeLeftEmu($Shape, $Emu) {
    $Shape.Left = ConvertTo-PtForEmu $Emu
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Every shape in the diagram moves straight down by 395923 EMU, except the
# title textbox ("Rectangle 29" / "High-level Representation") which also
# shifts slightly on X and moves down by a smaller amount.
$moves = @(
    @{ Name = "Cube 3";                       Y = 3440780 },
    @{ Name = "TextBox 5";                    Y = 2799758 },
    @{ Name = "Rectangle 6";                  Y = 2799758 },
    @{ Name = "Rectangle 10";                 Y = 2430426 },
    @{ Name = "Rectangle 12";                 Y = 2430426 },
    @{ Name = "Rectangle 14";                 Y = 2430426 },
    @{ Name = "Cube 19";                      Y = 3440780 },
    @{ Name = "Cube 20";                      Y = 3440780 },
    @{ Name = "Cube 22";                      Y = 3440780 },
    @{ Name = "Cube 23";                      Y = 3440780 },
    @{ Name = "Straight Arrow Connector 8";   Y = 3953363 },
    @{ Name = "Straight Arrow Connector 25";  Y = 3953363 },
    @{ Name = "Straight Arrow Connector 30";  Y = 3953363 },
    @{ Name = "Straight Arrow Connector 32";  Y = 3953363 },
    @{ Name = "Rectangle 1";                  Y = 2799758 },
    @{ Name = "Rectangle 17";                 Y = 2799758 },
    @{ Name = "Rectangle 2";                  Y = 2799758 },
    @{ Name = "Rectangle 21";                 Y = 2430426 },
    @{ Name = "TextBox 24";                   Y = 2430426 }
)

foreach ($move in $moves) {
    $shape = $s.Shapes.Item($move.Name)
    Set-ShapeTopEmu $shape $move.Y
}

# Title textbox: shifts on both axes.
$title = $s.Shapes.Item("Rectangle 29")
Set-ShapeLeftEmu $title 4061939
Set-ShapeTopEmu $title 1296972

# Best-effort: the canonical OOXML also gains an (empty) slide-guide list
# extension on the presentation itself (<p:extLst><p:ext uri="{EFAFB233-...}">
# <p15:sldGuideLst/></p:ext></p:extLst>), which PowerPoint writes once the
# Guides collection has been touched. Touch it here so it is emitted if the
# host supports it; this is inert (no guides are defined) either way.
try {
    $null = $p.Guides
} catch {
}
